{"js": "// Update the two-digit division worksheet: replace each division\n// problem's text with its new value (same cell, same formatting).\nconst replacements = [\n  [\"64\u00f79=\", \"42\u00f74=\"],\n  [\"98\u00f76=\", \"56\u00f75=\"],\n  [\"11\u00f76=\", \"14\u00f73=\"],\n  [\"39\u00f77=\", \"45\u00f77=\"],\n  [\"24\u00f72=\", \"88\u00f75=\"],\n  [\"52\u00f78=\", \"88\u00f72=\"],\n  [\"50\u00f79=\", \"48\u00f72=\"],\n  [\"64\u00f73=\", \"38\u00f75=\"],\n  [\"29\u00f73=\", \"62\u00f77=\"],\n  [\"88\u00f78=\", \"81\u00f78=\"],\n  [\"48\u00f76=\", \"10\u00f79=\"],\n  [\"70\u00f78=\", \"82\u00f74=\"],\n  [\"71\u00f72=\", \"77\u00f79=\"],\n  [\"53\u00f78=\", \"78\u00f77=\"],\n  [\"90\u00f76=\", \"19\u00f79=\"],\n  [\"87\u00f76=\", \"22\u00f72=\"],\n  [\"49\u00f75=\", \"39\u00f76=\"],\n  [\"67\u00f74=\", \"35\u00f74=\"],\n  [\"74\u00f78=\", \"66\u00f75=\"],\n  [\"56\u00f76=\", \"55\u00f73=\"],\n  [\"18\u00f73=\", \"17\u00f78=\"],\n  [\"86\u00f79=\", \"29\u00f76=\"],\n  [\"99\u00f78=\", \"72\u00f78=\"],\n  [\"58\u00f79=\", \"21\u00f78=\"],\n  [\"92\u00f79=\", \"50\u00f79=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the two-digit division worksheet: replace each division\n# problem's text with its new value (same cell, same formatting).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"64\u00f79=\", \"42\u00f74=\"),\n    @(\"98\u00f76=\", \"56\u00f75=\"),\n    @(\"11\u00f76=\", \"14\u00f73=\"),\n    @(\"39\u00f77=\", \"45\u00f77=\"),\n    @(\"24\u00f72=\", \"88\u00f75=\"),\n    @(\"52\u00f78=\", \"88\u00f72=\"),\n    @(\"50\u00f79=\", \"48\u00f72=\"),\n    @(\"64\u00f73=\", \"38\u00f75=\"),\n    @(\"29\u00f73=\", \"62\u00f77=\"),\n    @(\"88\u00f78=\", \"81\u00f78=\"),\n    @(\"48\u00f76=\", \"10\u00f79=\"),\n    @(\"70\u00f78=\", \"82\u00f74=\"),\n    @(\"71\u00f72=\", \"77\u00f79=\"),\n    @(\"53\u00f78=\", \"78\u00f77=\"),\n    @(\"90\u00f76=\", \"19\u00f79=\"),\n    @(\"87\u00f76=\", \"22\u00f72=\"),\n    @(\"49\u00f75=\", \"39\u00f76=\"),\n    @(\"67\u00f74=\", \"35\u00f74=\"),\n    @(\"74\u00f78=\", \"66\u00f75=\"),\n    @(\"56\u00f76=\", \"55\u00f73=\"),\n    @(\"18\u00f73=\", \"17\u00f78=\"),\n    @(\"86\u00f79=\", \"29\u00f76=\"),\n    @(\"99\u00f78=\", \"72\u00f78=\"),\n    @(\"58\u00f79=\", \"21\u00f78=\"),\n    @(\"92\u00f79=\", \"50\u00f79=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
